# The commit swaps the deck's two embedded themes: the "Integral" theme
# (ppt/theme/theme1.xml, wired to the slide master -> drives every slide)
# and the plain "Office Theme" (ppt/theme/theme2.xml, only wired to the
# notes master). After the edit theme1.xml carries the Office Theme's
# color scheme and theme2.xml carries what used to be theme1's (Integral).
#
# The PowerPoint object model only exposes write access to the *active*
# theme (the one backing the slide master / every slide) via the
# per-slide ThemeColorScheme color collection - there is no COM surface
# here that reaches the inert theme2.xml used solely by the notes master.
# So we reproduce the externally-visible part of the swap: push the
# "Office Theme" palette (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink)
# into the live color scheme, in clrScheme child order.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

$tcs.Item(1).RGB  = 0         # dk1      000000
$tcs.Item(2).RGB  = 16777215  # lt1      FFFFFF
$tcs.Item(3).RGB  = 6968388   # dk2      44546A
$tcs.Item(4).RGB  = 15132391  # lt2      E7E6E6
$tcs.Item(5).RGB  = 13998939  # accent1  5B9BD5
$tcs.Item(6).RGB  = 3243501   # accent2  ED7D31
$tcs.Item(7).RGB  = 10855845  # accent3  A5A5A5
$tcs.Item(8).RGB  = 49407     # accent4  FFC000
$tcs.Item(9).RGB  = 12874308  # accent5  4472C4
$tcs.Item(10).RGB = 4697456   # accent6  70AD47
$tcs.Item(11).RGB = 12673797  # hlink    0563C1
$tcs.Item(12).RGB = 7491477   # folHlink 954F72
